# Applies the "Anonymize fedcore" update:
#  - Adds top/bottom border to the first column of each merged metric
#    header group (C1-style), and top/bottom/right border to the last
#    column of each group (D1-style / G1-style).
#  - Renames the "fedcore" column header to "approach".
#  - Removes the stray empty cell G5 on the computational_comparison sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws2 = $wb.Worksheets.Item("computational_comparison")

# xlEdgeLeft=7, xlEdgeTop=8, xlEdgeBottom=9, xlEdgeRight=10
# xlContinuous=1, xlNone=-4142, xlThin=2 (weight)

function Set-TopBottomBorder($rng) {
    $rng.Style = "Normal"
    $rng.Borders.LineStyle = 1
    $rng.Borders.Weight = 2
    $rng.Borders.Item(7).LineStyle = -4142
    $rng.Borders.Item(10).LineStyle = -4142
}

function Set-TopBottomRightBorder($rng) {
    $rng.Style = "Normal"
    $rng.Borders.LineStyle = 1
    $rng.Borders.Weight = 2
    $rng.Borders.Item(7).LineStyle = -4142
}

# --- quality_comparison sheet ---
Set-TopBottomBorder $ws1.Range("C1")
Set-TopBottomRightBorder $ws1.Range("D1")
$ws1.Range("C2").Value = "approach"

# --- computational_comparison sheet ---
Set-TopBottomBorder $ws2.Range("C1")
Set-TopBottomRightBorder $ws2.Range("D1")
Set-TopBottomBorder $ws2.Range("F1")
Set-TopBottomRightBorder $ws2.Range("G1")

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty cell G5
$ws2.Range("G5").ClearContents()
